$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "jumping enemy " + "(spider?)" -> single run "jumping enemy (spider?)"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("jumping enemy")) {
        $full = $p.Range
        $splitOffset = $t.IndexOf("(spider?)")
        $tail = $d.Range($full.Start + $splitOffset, $full.Start + $splitOffset + 9)
        $tail.Delete()
        $joinPoint = $d.Range($full.Start + $splitOffset, $full.Start + $splitOffset)
        $joinPoint.InsertAfter("(spider?)")
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Add three new bullet paragraphs (numId 8 / ilvl 0) after the
#    "I suppose you don't know..." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("I suppose you don")) {
        $anchor = $p
        break
    }
}

$lsq = [string][char]0x2018
$rsq = [string][char]0x2019

# --- paragraph: 'Infinite' area that continually expands as it is explored?
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
$start = $newPara.Range.Start

$p1 = $d.Range($start, $start)
$p1.InsertAfter($lsq)
$r1 = $d.Range($start, $start + 1)
$r1.Font.Italic = $true

$p2 = $d.Range($start + 1, $start + 1)
$p2.InsertAfter("Infinite" + $rsq + " area that continually expands as it is explored?")
$r2End = $newPara.Range.End - 1
$r2 = $d.Range($start + 1, $r2End)
$r2.Font.Italic = $false

$r1b = $d.Range($start, $start + 1)
$r1b.Font.Italic = $false

# --- paragraph: Multi-part door gradually opened larger and larger?
$newPara.Range.InsertParagraphAfter()
$newPara2 = $newPara.Next()
$newPara2.Range.InsertAfter("Multi-part door gradually opened larger and larger?")

# --- paragraph: Grappling hook?
$newPara2.Range.InsertParagraphAfter()
$newPara3 = $newPara2.Next()
$newPara3.Range.InsertAfter("Grappling hook?")

# ---------------------------------------------------------------------------
# 3) Normal style: overflowPunct true -> false (ParagraphFormat.HangingPunctuation)
# ---------------------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $false
